$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) values for the 2020-09-04 data refresh.
# Values are stored as text in this dataset, so force text number format before assignment
# to avoid Excel auto-converting the numeric-looking strings into numbers.

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "195"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "458016.00"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1047"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3407744.33"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "426"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1776198.25"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "31"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "206643.82"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "8"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "50500.00"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "45"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "90000.00"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "61"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "159597.64"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "372"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1349761.71"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "37"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "190283.00"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "8"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42000.00"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "78"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "447894.00"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "30"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108621.84"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "53"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "306703.00"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "109"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "325135.17"
$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "619"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "2267961.21"
$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "272"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "1242878.76"
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "94"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = "568274.23"
$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "26"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "82220.65"
$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "130"
$ws.Range("D61").NumberFormat = "@"
$ws.Range("D61").Value = "898623.00"
$ws.Range("C70").NumberFormat = "@"
$ws.Range("C70").Value = "29"
$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = "64214.00"
$ws.Range("C71").NumberFormat = "@"
$ws.Range("C71").Value = "193"
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "477149.88"
$ws.Range("C72").NumberFormat = "@"
$ws.Range("C72").Value = "99"
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = "283926.60"
$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "15"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "30150.00"
$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "235"
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "603326.09"
$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "909"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "2935833.26"
$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "340"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "1387985.70"
$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "122"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "613484.52"
$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "35"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "76500.00"
$ws.Range("C95").NumberFormat = "@"
$ws.Range("C95").Value = "101"
$ws.Range("D95").NumberFormat = "@"
$ws.Range("D95").Value = "262578.00"
$ws.Range("C96").NumberFormat = "@"
$ws.Range("C96").Value = "420"
$ws.Range("D96").NumberFormat = "@"
$ws.Range("D96").Value = "1294526.53"
$ws.Range("C97").NumberFormat = "@"
$ws.Range("C97").Value = "179"
$ws.Range("D97").NumberFormat = "@"
$ws.Range("D97").Value = "734027.27"
$ws.Range("C98").NumberFormat = "@"
$ws.Range("C98").Value = "59"
$ws.Range("D98").NumberFormat = "@"
$ws.Range("D98").Value = "301911.73"
$ws.Range("C100").NumberFormat = "@"
$ws.Range("C100").Value = "15"
$ws.Range("D100").NumberFormat = "@"
$ws.Range("D100").Value = "31500.00"
